# Adds a new wave of data (28. 9. 2021) to both worksheets of the
# ZBP_12_obavy_ztrata_prace workbook:
#   - "data"   sheet gains column AH (percentages), rows 1-61
#   - "pocetR" sheet gains column AG (sample sizes),  rows 1-24
# and refreshes the "aktualizace" date stamped in the last row of each sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "data": new column AH
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("data")

# Header cell - copy formatting from the previous date header (AG1) so the
# new column matches the existing bold/centered/bordered style, then set
# the new wave's date label.
$wsData.Range("AG1").Copy($wsData.Range("AH1"))
$wsData.Range("AH1").Value = "28. 9. 2021"

# Percentage values for rows 2-61 (column AH), in row order.
$dataValues = @(
    0.54, 0.3, 0.16, 0.22, 0.28, 0.5, 0.57, 0.3, 0.13, 0.57,
    0.29, 0.14, 0.47, 0.31, 0.22, 0.58, 0.28, 0.14, 0.47, 0.31,
    0.22, 0.39, 0.35, 0.26, 0.48, 0.33, 0.19, 0.55, 0.29, 0.16,
    0.64, 0.25, 0.11, 0.47, 0.29, 0.24, 0.43, 0.35, 0.22, 0.61,
    0.28, 0.11, 0.67, 0.19, 0.14, 0.53, 0.31, 0.16, 0.62, 0.29,
    0.09, 0.47, 0.31, 0.22, 0.66, 0.24, 0.1, 0.71, 0.17, 0.12
)

$row = 2
foreach ($val in $dataValues) {
    $wsData.Cells.Item($row, 34).Value = $val   # column 34 = AH
    $row = $row + 1
}

# Refresh the "aktualizace" date mentioned in the trailing caption row (62).
$wsData.Range("A62").Value = "Život během pandemie, Obavy ze ztráty práce, % respondentů celkově a ve skupinách, aktualizace 6. 10. 2021"

# ---------------------------------------------------------------------
# Sheet "pocetR": new column AG
# ---------------------------------------------------------------------
$wsPocet = $wb.Worksheets.Item("pocetR")

# Header cell - copy formatting from the previous date header (AF1).
$wsPocet.Range("AF1").Copy($wsPocet.Range("AG1"))
$wsPocet.Range("AG1").Value = "28. 9. 2021"

# Sample-size values for rows 2-24 (column AG), in row order.
$pocetValues = @(
    1043, 88, 955, 814, 155, 7, 67, 782, 143, 59,
    59, 396, 403, 244, 114, 305, 357, 162, 290, 100,
    269, 144, 92
)

$row = 2
foreach ($val in $pocetValues) {
    $wsPocet.Cells.Item($row, 33).Value = $val   # column 33 = AG
    $row = $row + 1
}

# Refresh the "aktualizace" date mentioned in the trailing caption row (25).
$wsPocet.Range("A25").Value = "Život během pandemie, Obavy ze ztráty práce, velikost dotázaného souboru celkově a ve skupinách, aktualizace 6. 10. 2021"

# Row 25 carries a trailing run of blank (but present) cells through AF25;
# extend that run to the new AG25 cell by copying one of its blank
# neighbours so the used range/dimension includes AG25.
$wsPocet.Range("AF25").Copy($wsPocet.Range("AG25"))
